# Add a new "localdb" command type to the '#system' sheet.
#
# This inserts a new column N (pushing the existing N..AC block of
# command-type columns one column to the right, to O..AD), fills the new
# N column with the "localdb" header and its six commands, inserts a new
# entry ("localdb") into the sorted "target" list in column A (pushing
# rows 14..29 down to 15..30), and updates every defined name that
# pointed at a shifted column (plus adds the new "localdb" name).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("#system")

# ---------------------------------------------------------------------
# 1) Insert a new column before N. This shifts the existing N:AC data
#    (macro, mail, math... no wait just the per-command-type columns)
#    one column to the right, becoming O:AD, and leaves a blank column N
#    in their place, ready to be populated with the new "localdb" data.
# ---------------------------------------------------------------------
$ws.Columns("N:N").Insert()

# ---------------------------------------------------------------------
# 2) Populate the new column N with the "localdb" command type: a
#    header in N1 and its six commands in N2:N7.
# ---------------------------------------------------------------------
$ws.Range("N1").Value = "localdb"
$ws.Range("N2").Value = "cloneTable(var,source,target)"
$ws.Range("N3").Value = "dropTables(var,tables)"
$ws.Range("N4").Value = "exportCSV(sql,output)"
$ws.Range("N5").Value = "importRecords(var,sourceDb,sql,table)"
$ws.Range("N6").Value = "purge(var)"
$ws.Range("N7").Value = "runSQLs(var,sqls)"

# ---------------------------------------------------------------------
# 3) Insert "localdb" into the alphabetically-sorted list of command
#    types kept in column A (the "target" named range). "localdb" sorts
#    between "json" and "macro", i.e. at row 14, so rows 14..29 shift
#    down to 15..30 and the new entry lands on row 14.
# ---------------------------------------------------------------------
for ($r = 29; $r -ge 14; $r--) {
    $ws.Cells.Item($r + 1, 1).Value = $ws.Cells.Item($r, 1).Text
}
$ws.Cells.Item(14, 1).Value = "localdb"

# ---------------------------------------------------------------------
# 4) Fix up the defined names: every name that referred to one of the
#    columns which shifted right (N:AC -> O:AD) needs its reference
#    updated, "target" needs its row extended by one, "macro" needs to
#    point at its new column, and a brand-new "localdb" name is added.
# ---------------------------------------------------------------------
$wb.Names.Item("mail").RefersTo      = "='#system'!`$P`$2:`$P`$2"
$wb.Names.Item("number").RefersTo    = "='#system'!`$Q`$2:`$Q`$16"
$wb.Names.Item("pdf").RefersTo       = "='#system'!`$R`$2:`$R`$16"
$wb.Names.Item("rdbms").RefersTo     = "='#system'!`$S`$2:`$S`$7"
$wb.Names.Item("redis").RefersTo     = "='#system'!`$T`$2:`$T`$10"
$wb.Names.Item("sms").RefersTo       = "='#system'!`$U`$2:`$U`$2"
$wb.Names.Item("sound").RefersTo     = "='#system'!`$V`$2:`$V`$5"
$wb.Names.Item("ssh").RefersTo       = "='#system'!`$W`$2:`$W`$9"
$wb.Names.Item("step").RefersTo      = "='#system'!`$X`$2:`$X`$4"
$wb.Names.Item("target").RefersTo    = "='#system'!`$A`$2:`$A`$30"
$wb.Names.Item("web").RefersTo       = "='#system'!`$Y`$2:`$Y`$127"
$wb.Names.Item("webalert").RefersTo  = "='#system'!`$Z`$2:`$Z`$8"
$wb.Names.Item("webcookie").RefersTo = "='#system'!`$AA`$2:`$AA`$8"
$wb.Names.Item("ws").RefersTo        = "='#system'!`$AB`$2:`$AB`$17"
$wb.Names.Item("ws.async").RefersTo  = "='#system'!`$AC`$2:`$AC`$8"
$wb.Names.Item("xml").RefersTo       = "='#system'!`$AD`$2:`$AD`$21"
$wb.Names.Item("macro").RefersTo     = "='#system'!`$O`$2:`$O`$4"
$wb.Names.Add("localdb", "='#system'!`$N`$2:`$N`$7")
